$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update financial figures for rows 2-6 (years 2014-2018) to restated/corrected values
# Row 2
$ws.Range("D2").Value = 906
$ws.Range("E2").Value = 65
$ws.Range("F2").Value = 65
$ws.Range("G2").Value = 57
$ws.Range("H2").Value = 47
$ws.Range("I2").Value = 29
$ws.Range("J2").Value = 18
$ws.Range("K2").Value = 1140
$ws.Range("L2").Value = 646
$ws.Range("M2").Value = 494
$ws.Range("N2").Value = 374
$ws.Range("O2").Value = 120
$ws.Range("P2").Value = 53
$ws.Range("Q2").Value = 16
$ws.Range("R2").Value = -45
$ws.Range("S2").Value = -34
$ws.Range("T2").Value = 48
$ws.Range("U2").Value = -32
$ws.Range("V2").Value = 258
$ws.Range("W2").Value = 7.18
$ws.Range("X2").Value = 5.21
$ws.Range("Y2").Value = 8.07
$ws.Range("Z2").Value = 4.14
$ws.Range("AA2").Value = 130.64
$ws.Range("AB2").Value = 607.43
$ws.Range("AC2").Value = 2740
$ws.Range("AD2").Value = 7.59
$ws.Range("AE2").Value = 35432
$ws.Range("AF2").Value = 0.59
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 1.2
$ws.Range("AI2").Value = 9.12
$ws.Range("AJ2").Value = 1056000

# Row 3
$ws.Range("D3").Value = 1093
$ws.Range("E3").Value = 111
$ws.Range("F3").Value = 111
$ws.Range("G3").Value = 109
$ws.Range("H3").Value = 81
$ws.Range("I3").Value = 49
$ws.Range("J3").Value = 32
$ws.Range("K3").Value = 1142
$ws.Range("L3").Value = 576
$ws.Range("M3").Value = 566
$ws.Range("N3").Value = 419
$ws.Range("O3").Value = 147
$ws.Range("P3").Value = 53
$ws.Range("Q3").Value = 193
$ws.Range("R3").Value = -58
$ws.Range("S3").Value = -108
$ws.Range("T3").Value = 69
$ws.Range("U3").Value = 124
$ws.Range("V3").Value = 158
$ws.Range("W3").Value = 10.2
$ws.Range("X3").Value = 7.4
$ws.Range("Y3").Value = 12.26
$ws.Range("Z3").Value = 7.09
$ws.Range("AA3").Value = 101.87
$ws.Range("AB3").Value = 692.28
$ws.Range("AC3").Value = 4604
$ws.Range("AD3").Value = 7.25
$ws.Range("AE3").Value = 39680
$ws.Range("AF3").Value = 0.84
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 0.75
$ws.Range("AI3").Value = 5.43
$ws.Range("AJ3").Value = 1056000

# Row 4
$ws.Range("D4").Value = 1229
$ws.Range("E4").Value = 130
$ws.Range("F4").Value = 130
$ws.Range("G4").Value = 127
$ws.Range("H4").Value = 92
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 36
$ws.Range("K4").Value = 1226
$ws.Range("L4").Value = 576
$ws.Range("M4").Value = 650
$ws.Range("N4").Value = 473
$ws.Range("O4").Value = 177
$ws.Range("P4").Value = 53
$ws.Range("Q4").Value = 139
$ws.Range("R4").Value = -80
$ws.Range("S4").Value = -50
$ws.Range("T4").Value = 89
$ws.Range("U4").Value = 50
$ws.Range("V4").Value = 118
$ws.Range("W4").Value = 10.54
$ws.Range("X4").Value = 7.49
$ws.Range("Y4").Value = 12.62
$ws.Range("Z4").Value = 7.77
$ws.Range("AA4").Value = 88.64
$ws.Range("AB4").Value = 791.68
$ws.Range("AC4").Value = 5326
$ws.Range("AD4").Value = 7.43
$ws.Range("AE4").Value = 44757
$ws.Range("AF4").Value = 0.88
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 0.63
$ws.Range("AI4").Value = 4.69
$ws.Range("AJ4").Value = 1056000

# Row 5
$ws.Range("D5").Value = 1178
$ws.Range("E5").Value = 68
$ws.Range("F5").Value = 68
$ws.Range("G5").Value = 66
$ws.Range("H5").Value = 48
$ws.Range("I5").Value = 31
$ws.Range("J5").Value = 16
$ws.Range("K5").Value = 1254
$ws.Range("L5").Value = 565
$ws.Range("M5").Value = 689
$ws.Range("N5").Value = 500
$ws.Range("O5").Value = 189
$ws.Range("P5").Value = 53
$ws.Range("Q5").Value = -23
$ws.Range("R5").Value = -71
$ws.Range("S5").Value = 40
$ws.Range("T5").Value = 73
$ws.Range("U5").Value = -96
$ws.Range("V5").Value = 166
$ws.Range("W5").Value = 5.74
$ws.Range("X5").Value = 4.06
$ws.Range("Y5").Value = 6.43
$ws.Range("Z5").Value = 3.85
$ws.Range("AA5").Value = 82.05
$ws.Range("AB5").Value = 844.98
$ws.Range("AC5").Value = 2962
$ws.Range("AD5").Value = 9.99
$ws.Range("AE5").Value = 47358
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 0.84
$ws.Range("AI5").Value = 8.44
$ws.Range("AJ5").Value = 1056000

# Row 6
$ws.Range("D6").Value = 1194
$ws.Range("E6").Value = 45
$ws.Range("F6").Value = 45
$ws.Range("G6").Value = 47
$ws.Range("H6").Value = 36
$ws.Range("I6").Value = 30
$ws.Range("K6").Value = 1336
$ws.Range("L6").Value = 622
$ws.Range("M6").Value = 714
$ws.Range("N6").Value = 523
$ws.Range("P6").Value = 53
$ws.Range("Q6").Value = 98
$ws.Range("R6").Value = -45
$ws.Range("S6").Value = 13
$ws.Range("T6").Value = 55
$ws.Range("U6").Value = 42
$ws.Range("V6").Value = 184
$ws.Range("W6").Value = 3.81
$ws.Range("X6").Value = 3.05
$ws.Range("Y6").Value = 5.92
$ws.Range("Z6").Value = 2.82
$ws.Range("AA6").Value = 87.01
$ws.Range("AB6").Value = 889.63
$ws.Range("AC6").Value = 2868
$ws.Range("AD6").Value = 56.84
$ws.Range("AE6").Value = 49551
$ws.Range("AF6").Value = 3.29
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 0.15
$ws.Range("AI6").Value = 8.72
$ws.Range("AJ6").Value = 1056000

# Remove erroneous estimate-year data (rows 7-9 / 2019E-2021E) beyond column C
$ws.Range("D7:AJ9").ClearContents()

Write-Host "IFRS list corrected"